$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added contact-method columns.
$ws.Range("F1").Value = "Facebook"
$ws.Range("G1").Value = "LinkedIn"

# Example-row placeholder values, matching the existing "Twitter" example cell.
$ws.Range("F2").Value = "exampleperson"
$ws.Range("G2").Value = "exampleperson"

# Match the formatting Google Sheets already uses for the "Twitter" column
# (header + example row) so the two new columns look consistent with it.
$ws.Range("E1:E2").Copy()
$ws.Range("F1:G2").PasteSpecial(-4122)
